# Error Calculations and Plots
# Apply missing-data edits to the imputation worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Individual cell value corrections (rows 2-25) ---
$ws.Range("D3").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("F8").Value = 17.05
$ws.Range("F10").Value = 16.43
$ws.Range("F12").Value = ""
$ws.Range("F15").Value = 16.2
$ws.Range("F18").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("F25").Value = 16.6

# --- Remove two data rows entirely (RM 232 and SC 92), shifting the rest up ---
# Delete the lower row first so the higher row's index stays valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# --- Post-shift cell value corrections (now rows 26-33) ---
$ws.Range("C26").Value = 10.8
$ws.Range("C27").Value = ""
$ws.Range("F29").Value = ""
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
